# chore: update Sheets via scheduled runner
# Refreshes the Leve market-board profit figures (columns H-N) on each
# server sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) of the
# Table_<Sheet> listed leves, based on newly-pulled market data.

$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3380.0571
$ws.Range("I43").Value = 2306.6667
$ws.Range("K43").Value = 2306.6667
$ws.Range("M43").Value = -2237.6667
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H86").Value = 14954.728
$ws.Range("I86").Value = 14500.333
$ws.Range("J86").Value = 16999.5
$ws.Range("K86").Value = 14500.333
$ws.Range("L86").Value = 16999.5
$ws.Range("M86").Value = -13377.333
$ws.Range("N86").Value = -19245.5
$ws.Range("H89").Value = 14954.728
$ws.Range("I89").Value = 14500.333
$ws.Range("J89").Value = 16999.5
$ws.Range("K89").Value = 72501.66500000001
$ws.Range("L89").Value = 84997.5
$ws.Range("M89").Value = -66885.66500000001
$ws.Range("N89").Value = -96229.5
$ws.Range("H113").Value = 100001080
$ws.Range("I113").Value = 50000900
$ws.Range("J113").Value = 133334536
$ws.Range("K113").Value = 50000900
$ws.Range("L113").Value = 133334536
$ws.Range("M113").Value = -49997646
$ws.Range("N113").Value = -133341044
$ws.Range("H116").Value = 6806.8667
$ws.Range("I116").Value = 6309.7
$ws.Range("J116").Value = 7801.2
$ws.Range("K116").Value = 6309.7
$ws.Range("L116").Value = 7801.2
$ws.Range("M116").Value = -2867.7
$ws.Range("N116").Value = -14685.2
$ws.Range("H129").Value = 1045.091
$ws.Range("I129").Value = 660.5
$ws.Range("J129").Value = 2070.6667
$ws.Range("K129").Value = 1981.5
$ws.Range("L129").Value = 6212.000100000001
$ws.Range("M129").Value = 3018.5
$ws.Range("N129").Value = -16212.0001
$ws.Range("N137").Value = -53544
$ws.Range("H137").Value = 5136.4614
$ws.Range("I137").Value = 1833
$ws.Range("J137").Value = 16148
$ws.Range("K137").Value = 5499
$ws.Range("L137").Value = 48444
$ws.Range("M137").Value = -2949

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N2").ClearContents()
$ws.Range("H2").Value = 1589.5385
$ws.Range("I2").Value = 1589.5385
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1589.5385
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1476.5385
$ws.Range("H32").Value = 13159713
$ws.Range("I32").Value = 13515349
$ws.Range("K32").Value = 13515349
$ws.Range("M32").Value = -13515062
$ws.Range("H45").Value = 2478.5
$ws.Range("I45").Value = 2174.2666
$ws.Range("K45").Value = 2174.2666
$ws.Range("M45").Value = -1797.2666
$ws.Range("H61").Value = 41759828
$ws.Range("I61").Value = 83341000
$ws.Range("K61").Value = 83341000
$ws.Range("M61").Value = -83340788
$ws.Range("H74").Value = 14717686
$ws.Range("I74").Value = 31251774
$ws.Range("K74").Value = 31251774
$ws.Range("M74").Value = -31250900
$ws.Range("H77").Value = 14717686
$ws.Range("I77").Value = 31251774
$ws.Range("K77").Value = 156258870
$ws.Range("M77").Value = -156254502
$ws.Range("H101").Value = 209665
$ws.Range("J101").Value = 209665
$ws.Range("L101").Value = 209665
$ws.Range("N101").Value = -216155
$ws.Range("N116").ClearContents()
$ws.Range("H116").Value = 1589.5385
$ws.Range("I116").Value = 1589.5385
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1589.5385
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 704.4614999999999
$ws.Range("H122").Value = 1199.5
$ws.Range("I122").Value = 1199.5
$ws.Range("K122").Value = 3598.5
$ws.Range("M122").Value = -1148.5
$ws.Range("H132").Value = 8980
$ws.Range("I132").Value = 5611.6924
$ws.Range("K132").Value = 16835.0772
$ws.Range("M132").Value = -14305.0772
$ws.Range("H136").Value = 41759828
$ws.Range("I136").Value = 83341000
$ws.Range("K136").Value = 250023000
$ws.Range("M136").Value = -250020450

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N3").ClearContents()
$ws.Range("H3").Value = 1589.5385
$ws.Range("I3").Value = 1589.5385
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1589.5385
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1475.5385
$ws.Range("H5").Value = 3878.5
$ws.Range("I5").Value = 15000
$ws.Range("K5").Value = 15000
$ws.Range("M5").Value = -14887
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("H94").Value = 1976.7059
$ws.Range("I94").Value = 1993.375
$ws.Range("K94").Value = 1993.375
$ws.Range("M94").Value = -1542.375
$ws.Range("H96").Value = 38400.3
$ws.Range("J96").Value = 76901.25
$ws.Range("L96").Value = 76901.25
$ws.Range("N96").Value = -82393.25
$ws.Range("H107").Value = 1569.6428
$ws.Range("I107").Value = 1382.6923
$ws.Range("K107").Value = 1382.6923
$ws.Range("M107").Value = 537.3077000000001
$ws.Range("H134").Value = 60843.332
$ws.Range("I134").Value = 1735.9
$ws.Range("K134").Value = 5207.700000000001
$ws.Range("M134").Value = -2672.700000000001

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N58").Value = -2420
$ws.Range("H58").Value = 1507
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 2014
$ws.Range("K58").Value = 1000
$ws.Range("L58").Value = 2014
$ws.Range("M58").Value = -797
$ws.Range("H99").Value = 2706.8823
$ws.Range("I99").Value = 2085.1667
$ws.Range("J99").Value = 4199
$ws.Range("K99").Value = 2085.1667
$ws.Range("L99").Value = 4199
$ws.Range("M99").Value = -587.1667000000002
$ws.Range("N99").Value = -7195
$ws.Range("H126").Value = 2706.8823
$ws.Range("I126").Value = 2085.1667
$ws.Range("J126").Value = 4199
$ws.Range("K126").Value = 6255.500100000001
$ws.Range("L126").Value = 12597
$ws.Range("M126").Value = -3785.500100000001
$ws.Range("N126").Value = -17537
$ws.Range("H132").Value = 3563.1667
$ws.Range("I132").Value = 3563.1667
$ws.Range("K132").Value = 10689.5001
$ws.Range("M132").Value = -8159.500100000001
$ws.Range("N136").Value = -11142
$ws.Range("H136").Value = 1507
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 2014
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 6042
$ws.Range("M136").Value = -450

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 84999
$ws.Range("J37").Value = 84999
$ws.Range("L37").Value = 254997
$ws.Range("N37").Value = -255221
$ws.Range("H129").Value = 64642.266
$ws.Range("I129").Value = 722.25
$ws.Range("K129").Value = 2166.75
$ws.Range("M129").Value = 2833.25
$ws.Range("H131").Value = 5564.3887
$ws.Range("J131").Value = 3860.111
$ws.Range("L131").Value = 11580.333
$ws.Range("N131").Value = -21660.333
$ws.Range("H137").Value = 4984.3076
$ws.Range("I137").Value = 4828.4287
$ws.Range("J137").Value = 5166.1665
$ws.Range("K137").Value = 14485.2861
$ws.Range("L137").Value = 15498.4995
$ws.Range("M137").Value = -9385.286100000001
$ws.Range("N137").Value = -25698.4995

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2572.5264
$ws.Range("I80").Value = 2960.2727
$ws.Range("J80").Value = 2039.375
$ws.Range("K80").Value = 2960.2727
$ws.Range("L80").Value = 2039.375
$ws.Range("M80").Value = -1962.2727
$ws.Range("N80").Value = -4035.375
$ws.Range("H83").Value = 2572.5264
$ws.Range("I83").Value = 2960.2727
$ws.Range("J83").Value = 2039.375
$ws.Range("K83").Value = 14801.3635
$ws.Range("L83").Value = 10196.875
$ws.Range("M83").Value = -9809.363499999999
$ws.Range("N83").Value = -20180.875

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2165.0667
$ws.Range("I16").Value = 1999
$ws.Range("K16").Value = 1999
$ws.Range("M16").Value = -1829
$ws.Range("H46").Value = 4238.087
$ws.Range("J46").Value = 4605
$ws.Range("L46").Value = 4605
$ws.Range("N46").Value = -4981
$ws.Range("H74").Value = 45624.668
$ws.Range("J74").Value = 58437
$ws.Range("L74").Value = 58437
$ws.Range("N74").Value = -60433
$ws.Range("H77").Value = 45624.668
$ws.Range("J77").Value = 58437
$ws.Range("L77").Value = 175311
$ws.Range("N77").Value = -185295
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("H108").Value = 82000
$ws.Range("J108").Value = 82000
$ws.Range("L108").Value = 82000
$ws.Range("N108").Value = -89680
$ws.Range("H132").Value = 745384.8
$ws.Range("I132").Value = 1251814.5
$ws.Range("J132").Value = 166608
$ws.Range("K132").Value = 3755443.5
$ws.Range("L132").Value = 499824
$ws.Range("M132").Value = -3752913.5
$ws.Range("N132").Value = -504884
$ws.Range("H133").Value = 70775
$ws.Range("J133").Value = 70775
$ws.Range("L133").Value = 70775
$ws.Range("N133").Value = -75835

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N46").Value = -49962
$ws.Range("H46").Value = 49500
$ws.Range("J46").Value = 49500
$ws.Range("L46").Value = 49500
$ws.Range("H100").Value = 771.6
$ws.Range("I100").Value = 760.7059
$ws.Range("J100").Value = 833.3333
$ws.Range("K100").Value = 1521.4118
$ws.Range("L100").Value = 1666.6666
$ws.Range("M100").Value = -980.4118000000001
$ws.Range("N100").Value = -2748.6666
$ws.Range("H103").Value = 110367
$ws.Range("J103").Value = 110367
$ws.Range("L103").Value = 110367
$ws.Range("N103").Value = -112711
$ws.Range("N134").Value = -153570
$ws.Range("H134").Value = 49500
$ws.Range("J134").Value = 49500
$ws.Range("L134").Value = 148500
